$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted right before the current row 260,
# pushing that row (and every row after it) down by one. The newly opened
# row 260 gets the fresh data point.
$ws.Rows.Item(260).Insert()

$ws.Range("A260").Value2 = 1
$ws.Range("B260").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C260").Value2 = "Arica y Parinacota"
$ws.Range("D260").Value2 = 44841
$ws.Range("E260").Value2 = 15
$ws.Range("F260").Value2 = 100114013
$ws.Range("G260").Value2 = "Zanahoria"
$ws.Range("H260").Value2 = "Sin especificar"
$ws.Range("I260").Value2 = "Primera"
$ws.Range("J260").Value2 = 70
$ws.Range("K260").Value2 = 24000
$ws.Range("L260").Value2 = 25000
$ws.Range("M260").Value2 = 24500
$ws.Range("N260").Value2 = "`$/saco 20 kilos"
$ws.Range("O260").Value2 = "Valle de Camiña"
$ws.Range("P260").Value2 = 1225
$ws.Range("Q260").Value2 = 20
$ws.Range("R260").Value2 = "Hortaliza"
